# Update ticket/visitor counts on the "展览" (rId1) and "全部类型" (rId4) sheets.
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 8241
$wsExhibition.Range("F11").Value = 881

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 8241
$wsAllTypes.Range("F15").Value = 882
